$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Dark Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B2").Value = "Aggro %`n- 10"

$ws.Range("A3").Value = "Water Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B3").Value = "Ailment Resistance %`n5`nMagic Device only:Aggro %`n- 10"

$ws.Range("A4").Value = "Metal Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B4").Value = "Critical Rate`n5"

$ws.Range("A5").Value = "Wind Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B5").Value = "ASPD`n250`nKatana only:Critical Rate`n5"

$ws.Range("B6").Value = "Max HP %`n101`n-Handed Sword only:Fractional Barrier %`n10"

$ws.Range("A7").Value = "Fire Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B7").Value = "MATK %`n1`nStaff only:Magic Pierce %`n5"

$ws.Range("A8").Value = "Lightning Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B8").Value = "Stability %`n5`nKatana only:Accuracy %`n10"
